# Amelioration de l'affichage des statistiques
# Update the "Contenu du stage" table (rows 16-23, columns D/E/G) with the
# refreshed head-counts and recomputed percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell, far outside the used range, used to stage text values so that
# they can be pasted as literal text (avoids Excel's "smart" number/percent
# auto-conversion when assigning a string like "68 %" directly to .Value).
$staging = $ws.Range("ZZ100")

function Set-TextValue($range, [string]$text) {
    $staging.Value = '="' + $text + '"'
    $staging.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# C# : 19 -> 17 students, 76 % -> 68 %
$ws.Range("E16").Value = 17
Set-TextValue $ws.Range("G16") "68 %"

# COBOL : unchanged (1 student, 4 %)

# C++ : unchanged (0 students, 0 %)

# ASSEMBLEUR : 3 -> 4 students, 12 % -> 16 %
$ws.Range("E19").Value = 4
Set-TextValue $ws.Range("G19") "16 %"

# ANDROID : 0 -> 2 students, 0 % -> 8 %
$ws.Range("E20").Value = 2
Set-TextValue $ws.Range("G20") "8 %"

# JEE : unchanged (1 student, 4 %)

# DELPHI : unchanged (0 students, 0 %)

# PHP5 : 1 -> 0 students, 4 % -> 0 %
$ws.Range("E23").Value = 0
Set-TextValue $ws.Range("G23") "0 %"

$staging.ClearContents()

$wb.Save()
